$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 908
$ws.Range("F3").Value = 1479
$ws.Range("G3").Value = 45
$ws.Range("F4").Value = 1146
$ws.Range("G4").Value = 70
$ws.Range("F6").Value = 231
$ws.Range("F7").Value = 8
$ws.Range("F8").Value = 697
$ws.Range("F9").Value = 288
$ws.Range("F11").Value = 107
$ws.Range("F13").Value = 168
$ws.Range("F14").Value = 3671
$ws.Range("F15").Value = 22
$ws.Range("F19").Value = 516
$ws.Range("F24").Value = 685
$ws.Range("F25").Value = 67
$ws.Range("F29").Value = 1633
$ws.Range("F30").Value = 365

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 33
$ws.Range("F7").Value = 245

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 396
$ws.Range("F3").Value = 121

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 396
$ws.Range("F3").Value = 908
$ws.Range("F4").Value = 1479
$ws.Range("G4").Value = 45
$ws.Range("F5").Value = 1146
$ws.Range("G5").Value = 70
$ws.Range("F8").Value = 121
$ws.Range("F10").Value = 231
$ws.Range("F11").Value = 8
$ws.Range("F12").Value = 697
$ws.Range("F14").Value = 288
$ws.Range("F16").Value = 107
$ws.Range("F18").Value = 168
$ws.Range("F19").Value = 3671
$ws.Range("F20").Value = 22
$ws.Range("F25").Value = 516
$ws.Range("F28").Value = 33
$ws.Range("F31").Value = 245
$ws.Range("F34").Value = 685
$ws.Range("F38").Value = 67
$ws.Range("F42").Value = 1633
$ws.Range("F43").Value = 365
